# Insert two new rows at 384-385 (existing rows 384.. shift down to 386..),
# then populate the two new rows with the new record data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("384:385").Insert()

# New row 384
$ws.Range("A384").Value = 4
$ws.Range("B384").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C384").Value = "Los Lagos"
$ws.Range("D384").Value = 44932
$ws.Range("E384").Value = 10
$ws.Range("F384").Value = 100114013
$ws.Range("G384").Value = "Zanahoria"
$ws.Range("H384").Value = "Sin especificar"
$ws.Range("I384").Value = "Primera"
$ws.Range("J384").Value = 450
$ws.Range("K384").Value = 14000
$ws.Range("L384").Value = 14000
$ws.Range("M384").Value = 14000
$ws.Range("N384").Value = "$/saco 20 kilos"
$ws.Range("O384").Value = "Región Metropolitana"
$ws.Range("P384").Value = 700
$ws.Range("Q384").Value = 20
$ws.Range("R384").Value = "Hortaliza"

# New row 385
$ws.Range("A385").Value = 4
$ws.Range("B385").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C385").Value = "Los Lagos"
$ws.Range("D385").Value = 44932
$ws.Range("E385").Value = 10
$ws.Range("F385").Value = 100114013
$ws.Range("G385").Value = "Zanahoria"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 450
$ws.Range("K385").Value = 15000
$ws.Range("L385").Value = 15000
$ws.Range("M385").Value = 15000
$ws.Range("N385").Value = "$/saco 20 kilos"
$ws.Range("O385").Value = "Región de Coquimbo"
$ws.Range("P385").Value = 750
$ws.Range("Q385").Value = 20
$ws.Range("R385").Value = "Hortaliza"
